$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells keep their original text formatting
# (values like "1.003" or "22.50" would otherwise be auto-converted to numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.075.05"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.929.53"
$ws.Range("E3").Value = "  +2.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.83"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4604"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07743"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9833"
$ws.Range("E10").Value = "  +2.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.50"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.921.84"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.979"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.696"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07025"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.42"
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009531"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.76"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.090.78"
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.352"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.97"
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.078"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.90"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.12"
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.701"
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "118.13"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.857"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09337"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8672"
$ws.Range("E31").Value = "  +2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.128"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.255"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.047"
$ws.Range("E34").Value = "  -0.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05720"
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02050"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.051"
$ws.Range("E39").Value = "  +13.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.551"
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5529"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1757"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000002974"
$ws.Range("E43").Value = "  +2.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.403"
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.243"
$ws.Range("E45").Value = "  +8.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5201"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.25"
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06923"
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.785"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.61"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.12"
$ws.Range("E51").Value = "  +1.14%  "
